$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 197; this pushes the existing rows
# 197..246 down to 198..247, preserving all of their data/formatting.
$ws.Rows(197).EntireRow.Insert()

# Populate the newly inserted row 197 with the new record.
$ws.Range("A197").Value = 10
$ws.Range("B197").Value = "Vega Modelo de Temuco"
$ws.Range("C197").Value = "La Araucanía"
$ws.Range("D197").Value = 44988
$ws.Range("E197").Value = 9
$ws.Range("F197").Value = 100112012
$ws.Range("G197").Value = "Espinaca"
$ws.Range("H197").Value = "Sin especificar"
$ws.Range("I197").Value = "Primera"
$ws.Range("J197").Value = 20
$ws.Range("K197").Value = 12000
$ws.Range("L197").Value = 12000
$ws.Range("M197").Value = 12000
$ws.Range("N197").Value = "$/docena de atados"
$ws.Range("O197").Value = "Región de La Araucanía"
$ws.Range("P197").Value = 4000
$ws.Range("Q197").Value = 3
$ws.Range("R197").Value = "Hortaliza"
